$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    100 = 2155269.754
    101 = 546524.8496
    102 = 244173.8486
    103 = 1100026.488
    104 = 1657899.811
    105 = 4852606.837
    106 = 26859.12335
    107 = 1920085.387
    114 = 3483.011366
    115 = 8314196.243
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $range = $ws.Range("J$row`:AS$row")
    $range.Value = $value
}
